$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28: The Writing Is Not on the Wall
$ws.Range("H28").Value = 16288.77
$ws.Range("I28").Value = 25919
$ws.Range("K28").Value = 25919
$ws.Range("M28").Value = -25434
# Row 41: The Write Stuff
$ws.Range("H41").Value = 567
$ws.Range("I41").Value = 146.33333
$ws.Range("J41").Value = 882.5
$ws.Range("K41").Value = 146.33333
$ws.Range("L41").Value = 882.5
$ws.Range("M41").Value = 293.66667
$ws.Range("N41").Value = -1762.5
# Row 86: Filling in the Blanks
$ws.Range("H86").Value = 2989.4285
$ws.Range("I86").Value = 2734.4443
$ws.Range("J86").Value = 3448.4
$ws.Range("K86").Value = 2734.4443
$ws.Range("L86").Value = 3448.4
$ws.Range("M86").Value = -1611.4443
$ws.Range("N86").Value = -5694.4
# Row 89: Ink into Antiquity (L)
$ws.Range("H89").Value = 2989.4285
$ws.Range("I89").Value = 2734.4443
$ws.Range("J89").Value = 3448.4
$ws.Range("K89").Value = 13672.2215
$ws.Range("L89").Value = 17242
$ws.Range("M89").Value = -8056.2215
$ws.Range("N89").Value = -28474
# Row 92: Whinier than the Sword
$ws.Range("H92").Value = 1085.1904
$ws.Range("I92").Value = 1138.9412
$ws.Range("J92").Value = 856.75
$ws.Range("K92").Value = 1138.9412
$ws.Range("L92").Value = 856.75
$ws.Range("M92").Value = 109.0588
$ws.Range("N92").Value = -3352.75
# Row 98: The Dotted Line
$ws.Range("H98").Value = 10178.444
$ws.Range("I98").Value = 11900.4
$ws.Range("K98").Value = 11900.4
$ws.Range("M98").Value = -10402.4
# Row 111: An Eye for Healing
$ws.Range("H111").Value = 111111790
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 111111790
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 333335370
$ws.Range("N111").Value = -333341504
$ws.Range("M111").ClearContents()
# Row 122: Wishful Inking
$ws.Range("H122").Value = 10178.444
$ws.Range("I122").Value = 11900.4
$ws.Range("K122").Value = 35701.2
$ws.Range("M122").Value = -33251.2
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 1470.2391
$ws.Range("I132").Value = 1081.2439
$ws.Range("J132").Value = 4660
$ws.Range("K132").Value = 3243.7317
$ws.Range("L132").Value = 13980
$ws.Range("M132").Value = -713.7316999999998
$ws.Range("N132").Value = -19040

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 26: Night Squawker
$ws.Range("H26").Value = 1013.2857
$ws.Range("I26").Value = 1013.2857
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 1013.2857
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -683.2857
$ws.Range("N26").ClearContents()
# Row 122: Haste for High Durium
$ws.Range("H122").Value = 5000
$ws.Range("J122").Value = 5000
$ws.Range("L122").Value = 15000
$ws.Range("N122").Value = -19900
# Row 135: Forgiveness for My Shins
$ws.Range("H135").Value = 84809.664
$ws.Range("J135").Value = 84809.664
$ws.Range("L135").Value = 84809.664
$ws.Range("N135").Value = -94949.664

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 36: I Saw What You Did There
$ws.Range("H36").Value = 1000
$ws.Range("I36").Value = 1000
$ws.Range("K36").Value = 1000
$ws.Range("M36").Value = -466
# Row 62: Barring the Gates to Foundation
$ws.Range("H62").Value = 42857
$ws.Range("I62").Value = 39999
$ws.Range("J62").Value = 44000.2
$ws.Range("K62").Value = 39999
$ws.Range("L62").Value = 44000.2
$ws.Range("M62").Value = -39313
$ws.Range("N62").Value = -45372.2
# Row 65: Starting Young (L)
$ws.Range("H65").Value = 42857
$ws.Range("I65").Value = 39999
$ws.Range("J65").Value = 44000.2
$ws.Range("K65").Value = 119997
$ws.Range("L65").Value = 132000.6
$ws.Range("M65").Value = -116565
$ws.Range("N65").Value = -138864.6
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 2125589
$ws.Range("I86").Value = 18343.143
$ws.Range("K86").Value = 18343.143
$ws.Range("M86").Value = -17220.143
# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 2125589
$ws.Range("I89").Value = 18343.143
$ws.Range("K89").Value = 91715.715
$ws.Range("M89").Value = -86099.715
# Row 123: Archon Denied
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 21: Nightmare on My Street
$ws.Range("H21").Value = 2971.6667
$ws.Range("J21").Value = 2971.6667
$ws.Range("L21").Value = 2971.6667
$ws.Range("N21").Value = -3441.6667
# Row 31: Wall Not Found
$ws.Range("H31").Value = 7287.2144
$ws.Range("I31").Value = 2753.125
$ws.Range("K31").Value = 2753.125
$ws.Range("M31").Value = -2458.125
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 7287.2144
$ws.Range("I34").Value = 2753.125
$ws.Range("K34").Value = 2753.125
$ws.Range("M34").Value = -2551.125
# Row 55: Ready for a Rematch
$ws.Range("H55").Value = 15499.5
$ws.Range("J55").Value = 15499.5
$ws.Range("L55").Value = 15499.5
$ws.Range("N55").Value = -16129.5
# Row 68: Do You Even String Bow
$ws.Range("H68").Value = 63749.75
$ws.Range("I68").Value = 59999
$ws.Range("K68").Value = 59999
$ws.Range("M68").Value = -59250
# Row 71: Win One Bow, Get Three Free (L)
$ws.Range("H71").Value = 63749.75
$ws.Range("I71").Value = 59999
$ws.Range("K71").Value = 179997
$ws.Range("M71").Value = -176253
# Row 99: O Pine
$ws.Range("H99").Value = 2386.0833
$ws.Range("I99").Value = 2331.4546
$ws.Range("J99").Value = 2987
$ws.Range("K99").Value = 2331.4546
$ws.Range("L99").Value = 2987
$ws.Range("M99").Value = -833.4546
$ws.Range("N99").Value = -5983
# Row 103: Spare a Rod and Spoil the Fishers
$ws.Range("H103").Value = 9503.666999999999
$ws.Range("I103").Value = 9503.666999999999
$ws.Range("K103").Value = 9503.666999999999
$ws.Range("M103").Value = -8331.666999999999
# Row 126: A Better Conductor
$ws.Range("H126").Value = 2386.0833
$ws.Range("I126").Value = 2331.4546
$ws.Range("J126").Value = 2987
$ws.Range("K126").Value = 6994.3638
$ws.Range("L126").Value = 8961
$ws.Range("M126").Value = -4524.3638
$ws.Range("N126").Value = -13901
# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 1484.3112
$ws.Range("I132").Value = 1480.1708
$ws.Range("K132").Value = 4440.512400000001
$ws.Range("M132").Value = -1910.512400000001
# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 4032.9167
$ws.Range("I134").Value = 3490.4546
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 10471.3638
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -7936.363799999999
$ws.Range("N134").Value = -35070

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 23: Sweet Smell of Success
$ws.Range("H23").Value = 4242.6
$ws.Range("I23").Value = 5055.3335
$ws.Range("J23").Value = 3700.7778
$ws.Range("K23").Value = 15166.0005
$ws.Range("L23").Value = 11102.3334
$ws.Range("M23").Value = -14931.0005
$ws.Range("N23").Value = -11572.3334
# Row 137: Creative Chocolate
$ws.Range("H137").Value = 4759.222
$ws.Range("I137").Value = 3750
$ws.Range("J137").Value = 4839.96
$ws.Range("K137").Value = 11250
$ws.Range("L137").Value = 14519.88
$ws.Range("M137").Value = -6150
$ws.Range("N137").Value = -24719.88

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 46: Burning the Midnight Oil
$ws.Range("H46").Value = 33408.2
$ws.Range("J46").Value = 43333.332
$ws.Range("L46").Value = 43333.332
$ws.Range("N46").Value = -43645.332
# Row 53: North Ore South
$ws.Range("H53").Value = 24000
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
# Row 57: Gold Is So Last Year
$ws.Range("H57").Value = 10000000
$ws.Range("J57").Value = 10000000
$ws.Range("L57").Value = 10000000
$ws.Range("N57").Value = -10001640
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 5688.1875
$ws.Range("I80").Value = 3216.1428
$ws.Range("K80").Value = 3216.1428
$ws.Range("M80").Value = -2218.1428
# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 5688.1875
$ws.Range("I83").Value = 3216.1428
$ws.Range("K83").Value = 16080.714
$ws.Range("M83").Value = -11088.714
# Row 132: On Board for Lar
$ws.Range("H132").Value = 4105.2856
$ws.Range("I132").Value = 4105.2856
$ws.Range("K132").Value = 12315.8568
$ws.Range("M132").Value = -9785.856800000001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 4: Sole Traders
$ws.Range("H4").Value = 7009
$ws.Range("I4").Value = 7009
$ws.Range("K4").Value = 7009
$ws.Range("M4").Value = -6896
# Row 28: My Sole to Take
$ws.Range("H28").Value = 7009
$ws.Range("I28").Value = 7009
$ws.Range("K28").Value = 7009
$ws.Range("M28").Value = -6777
# Row 37: Quicker than Sand
$ws.Range("H37").Value = 7009
$ws.Range("I37").Value = 7009
$ws.Range("K37").Value = 7009
$ws.Range("M37").Value = -6902
# Row 43: Subordinate Clause
$ws.Range("H43").Value = 23342.229
$ws.Range("J43").Value = 24117.295
$ws.Range("L43").Value = 24117.295
$ws.Range("N43").Value = -24503.295
# Row 46: Supply Side Logic
$ws.Range("H46").Value = 2907.2307
$ws.Range("I46").Value = 1641.6666
$ws.Range("J46").Value = 3992
$ws.Range("K46").Value = 1641.6666
$ws.Range("L46").Value = 3992
$ws.Range("M46").Value = -1453.6666
$ws.Range("N46").Value = -4368
# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 3513.3794
$ws.Range("I61").Value = 2599.4443
$ws.Range("J61").Value = 5008.909
$ws.Range("K61").Value = 2599.4443
$ws.Range("L61").Value = 5008.909
$ws.Range("M61").Value = -2397.4443
$ws.Range("N61").Value = -5412.909
# Row 113: Peace in Rest
$ws.Range("H113").Value = 3513.3794
$ws.Range("I113").Value = 2599.4443
$ws.Range("J113").Value = 5008.909
$ws.Range("K113").Value = 2599.4443
$ws.Range("L113").Value = 5008.909
$ws.Range("M113").Value = -429.4443000000001
$ws.Range("N113").Value = -9348.909
# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 5143.8125
$ws.Range("I132").Value = 4691.8335
$ws.Range("K132").Value = 14075.5005
$ws.Range("M132").Value = -11545.5005

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 6773.423
$ws.Range("I132").Value = 6773.423
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 20320.269
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -17790.269
$ws.Range("N132").ClearContents()
# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 1019.8333
$ws.Range("I136").Value = 1030
$ws.Range("J136").Value = 999.5
$ws.Range("K136").Value = 3090
$ws.Range("L136").Value = 2998.5
$ws.Range("M136").Value = -540
$ws.Range("N136").Value = -8098.5
